# Auto-generated edit script: updates cryptos table rows 2-51
# (price/volume refresh + one-row insertion of "Aave", removal of "OKB")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.062.90"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.884.09"
$ws.Range("E3").Value = "  +0.76%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.9974"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5: BNB
$ws.Range("D5").Value = "'243.96"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6: USDC
$ws.Range("D6").Value = "'0.9969"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.4948"
$ws.Range("E7").Value = "  -0.43%  "

# Row 8: Cardano
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2921"
$ws.Range("E8").Value = "  +3.08%  "

# Row 9: Dogecoin
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06636"
$ws.Range("E9").Value = "  +1.44%  "

# Row 10: WrappedEther
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.880.92"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11: Solana
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'17.03"
$ws.Range("E11").Value = "  +0.34%  "

# Row 12: TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07203"
$ws.Range("E12").Value = "  -0.05%  "

# Row 13: Polygon
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.6657"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14: Litecoin
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'85.67"
$ws.Range("E14").Value = "  +1.07%  "

# Row 15: Polkadot
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.873"
$ws.Range("E15").Value = "  +1.58%  "

# Row 16: WrappedBTC
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "30.036.55"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17: ShibaInu
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000007880"
$ws.Range("E17").Value = "  +5.24%  "

# Row 18: Dai
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9967"
$ws.Range("E18").Value = "  -0.41%  "

# Row 19: Avalanche
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'12.83"
$ws.Range("E19").Value = "  +0.19%  "

# Row 20: WrappedliquidstakedEther2.0
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.123.11"
$ws.Range("E20").Value = "  +0.55%  "

# Row 21: BinanceUSD
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "'0.9969"
$ws.Range("E21").Value = "  -0.24%  "

# Row 22: Uniswap
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.774"
$ws.Range("E22").Value = "  +0.79%  "

# Row 23: Cosmos
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'9.184"
$ws.Range("E23").Value = "  +1.67%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "'5.607"
$ws.Range("E24").Value = "  +2.07%  "

# Row 25: Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'149.45"
$ws.Range("E25").Value = "  +3.21%  "

# Row 26: BitcoinCash
$ws.Range("B26").Value = "BitcoinCash"
$ws.Range("C26").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D26").Value = "'135.64"
$ws.Range("E26").Value = "  +0.17%  "

# Row 27: EthereumClassic
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'16.82"
$ws.Range("E27").Value = "  +0.47%  "

# Row 28: LidoDAOToken
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -1.36%  "

# Row 29: Toncoin
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.383"
$ws.Range("E29").Value = "  -0.64%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.199"
$ws.Range("E30").Value = "  -0.82%  "

# Row 31: Stellar
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08659"
$ws.Range("E31").Value = "  +0.44%  "

# Row 32: Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.961"
$ws.Range("E32").Value = "  +1.92%  "

# Row 33: Hedera
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04992"
$ws.Range("E33").Value = "  -1.52%  "

# Row 34: ARBITRUM
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.113"
$ws.Range("E34").Value = "  -1.28%  "

# Row 35: ImmutableX
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7027"
$ws.Range("E35").Value = "  +2.73%  "

# Row 36: HuobiToken
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.654"
$ws.Range("E36").Value = "  -1.36%  "

# Row 37: RenderToken
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.217"
$ws.Range("E37").Value = "  -4.71%  "

# Row 38: MXToken
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.688"
$ws.Range("E38").Value = "  -1.75%  "

# Row 39: TrustWalletToken
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.9361"
$ws.Range("E39").Value = "  -2.77%  "

# Row 40: VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01646"
$ws.Range("E40").Value = "  +1.08%  "

# Row 41: FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.992"
$ws.Range("E41").Value = "  -1.76%  "

# Row 42: PaxDollar
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9980"
$ws.Range("E42").Value = "  -0.20%  "

# Row 43: TheSandbox
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4208"
$ws.Range("E43").Value = "  +0.59%  "

# Row 44: Quant
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'101.82"
$ws.Range("E44").Value = "  -1.80%  "

# Row 45: Aptos
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'7.603"
$ws.Range("E45").Value = "  +2.42%  "

# Row 46: Algorand
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1264"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47: Cronos
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05716"
$ws.Range("E47").Value = "  +1.75%  "

# Row 48: Elrond
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'32.55"
$ws.Range("E48").Value = "  +0.44%  "

# Row 49: EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.230"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50: Aave
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'56.04"
$ws.Range("E50").Value = "  +2.45%  "

# Row 51: Decentraland
$ws.Range("D51").Value = "'0.3717"
$ws.Range("E51").Value = "  -0.28%  "
